$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$values = @(
    "ABCG2",
    "ACSL6",
    "AHSP",
    "AQP9",
    "ARG1",
    "BLVRB",
    "BRCA2",
    "BSG",
    "C1R",
    "CD209",
    "CD9",
    "CLEC5A",
    "CROT",
    "EGF",
    "FECH",
    "FSCN1",
    "FZD5",
    "G0S2",
    "GPX2",
    "GSTM3",
    "HBD",
    "HEBP1",
    "HMBS",
    "ICAM1",
    "INHBB",
    "JAM3",
    "MAFF",
    "MARCH8",
    "MEIS1",
    "MYL4",
    "MYL9",
    "PEX7",
    "PTP4A3",
    "RAD51C",
    "RBM38",
    "RETN",
    "RPL14",
    "SGCD",
    "SLC22A4",
    "SSR3",
    "TIMP1",
    "TMEM176A",
    "TPM1",
    "TRIM10",
    "TRIM25",
    "TSC22D1",
    "TSPAN4",
    "DACT1",
    "DBNDD1",
    "ERC1",
    "GPR25",
    "NLRP2",
    "RNASE6",
    "SLC14A1",
    "TTTY15",
    "AGAP1",
    "ALPL",
    "APOBEC3G",
    "BTN3A2",
    "CD79B",
    "ECHDC3",
    "ERO1B",
    "FAM13A",
    "GALNT6",
    "GPM6A",
    "HLA-DPA1",
    "HNRNPH1",
    "IL12RB2",
    "ITPR1",
    "KLHDC4",
    "LGALSL",
    "MRC2",
    "NEBL",
    "POLR1D",
    "PRKAR2B",
    "RUNX3",
    "SLC7A8",
    "SNORA21",
    "SPATA20",
    "STRN3",
    "TPP2",
    "ACAA2",
    "ACP1",
    "ACP2",
    "ADM",
    "ALB",
    "ASNS",
    "BPGM",
    "CACNA2D2",
    "CALD1",
    "CCND1",
    "CD28",
    "CDC25B",
    "DGAT1",
    "DHFR",
    "DHX16",
    "DLG5",
    "EBP",
    "ENDOD1",
    "ENO3",
    "ETV1",
    "GATA1",
    "GDE1",
    "GGA2",
    "GNG11",
    "HBQ1",
    "HDGF",
    "HLA-F",
    "IRF5",
    "KLHDC3",
    "LAMP2",
    "MARK3",
    "MCM5",
    "MKRN1",
    "MMP1",
    "MPP1",
    "MPZL2",
    "OPN3",
    "PCDH7",
    "PDK3",
    "PGRMC1",
    "PLOD2",
    "PMPCA",
    "PPP2R5B",
    "PQLC1",
    "PREB",
    "RIOK3",
    "SMTN",
    "SNCA",
    "SNRPA",
    "THBD",
    "TMCC2",
    "TMEM158",
    "TSPAN5",
    "XYLT2",
    "BEX3",
    "CLEC10A",
    "CNTN5",
    "ENY2",
    "FHL2",
    "LAMTOR2",
    "DOC2A",
    "MARC1",
    "PASK",
    "PLCL1",
    "PRR5L",
    "TNFRSF21",
    "ACADVL",
    "CSF1",
    "DUT",
    "IDH1",
    "MCL1",
    "PGK1",
    "PNP",
    "RB1",
    "SLC6A8",
    "SPTBN2",
    "TAL1",
    "TCF7L2",
    "TFDP2",
    "TGIF2",
    "UQCR10",
    "VEGFA",
    "ATP2A3",
    "GZMH",
    "MAPKAPK2",
    "RGS9",
    "ZSWIM8"
)

$n = $values.Length
$range = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($n, 1))

# Build a 2D array (n rows x 1 col) for a single bulk write
$arr = New-Object "object[,]" $n,1
for ($i = 0; $i -lt $n; $i++) {
    $arr[$i, 0] = $values[$i]
}

$range.Value = $arr

$wb.Save()
